$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 2835675
$ws.Cells.Item(2, 2).Value = 396830
$ws.Cells.Item(2, 3).Value = 516517
$ws.Cells.Item(3, 1).Value = 690706
$ws.Cells.Item(3, 2).Value = 31272
$ws.Cells.Item(3, 3).Value = 54819
$ws.Cells.Item(4, 1).Value = 513873
$ws.Cells.Item(4, 2).Value = 18858
$ws.Cells.Item(4, 3).Value = 70746
$ws.Cells.Item(5, 1).Value = 62932
$ws.Cells.Item(5, 2).Value = 29557
$ws.Cells.Item(5, 3).Value = 66616
$ws.Cells.Item(6, 1).Value = 530932
$ws.Cells.Item(6, 2).Value = 12548
$ws.Cells.Item(6, 3).Value = 38633
$ws.Cells.Item(7, 1).Value = 3368250
$ws.Cells.Item(7, 2).Value = 24746
$ws.Cells.Item(7, 3).Value = 58113
$ws.Cells.Item(8, 1).Value = 44191
$ws.Cells.Item(8, 2).Value = 35852
$ws.Cells.Item(8, 3).Value = 49424
$ws.Cells.Item(9, 1).Value = 4024608
$ws.Cells.Item(9, 2).Value = 37692
$ws.Cells.Item(9, 3).Value = 58406
$ws.Cells.Item(10, 1).Value = 260116
$ws.Cells.Item(10, 2).Value = 15265
$ws.Cells.Item(10, 3).Value = 28442
$ws.Cells.Item(11, 1).Value = 477943
$ws.Cells.Item(11, 2).Value = 7259
$ws.Cells.Item(11, 3).Value = 24126
$ws.Cells.Item(12, 1).Value = 1503173
$ws.Cells.Item(12, 2).Value = 9984
$ws.Cells.Item(12, 3).Value = 292995
$ws.Cells.Item(13, 1).Value = 420217
$ws.Cells.Item(13, 2).Value = 13868
$ws.Cells.Item(13, 3).Value = 31842
$ws.Cells.Item(14, 1).Value = 175397
$ws.Cells.Item(14, 2).Value = 9468
$ws.Cells.Item(14, 3).Value = 37975
$ws.Cells.Item(15, 1).Value = 44574
$ws.Cells.Item(15, 2).Value = 15145
$ws.Cells.Item(15, 3).Value = 32382
$ws.Cells.Item(16, 1).Value = 440140
$ws.Cells.Item(16, 2).Value = 9746
$ws.Cells.Item(16, 3).Value = 36070
$ws.Cells.Item(17, 1).Value = 44162
$ws.Cells.Item(17, 2).Value = 9466
$ws.Cells.Item(17, 3).Value = 38740
$ws.Cells.Item(18, 1).Value = 52174
$ws.Cells.Item(18, 2).Value = 12546
$ws.Cells.Item(18, 3).Value = 31203
$ws.Cells.Item(19, 1).Value = 419443
$ws.Cells.Item(19, 2).Value = 13490
$ws.Cells.Item(19, 3).Value = 34802
$ws.Cells.Item(20, 1).Value = 253996
$ws.Cells.Item(20, 2).Value = 13379
$ws.Cells.Item(20, 3).Value = 32303
$ws.Cells.Item(21, 1).Value = 402142
$ws.Cells.Item(21, 2).Value = 4687
$ws.Cells.Item(21, 3).Value = 31712
$ws.Cells.Item(22, 1).Value = 398354
$ws.Cells.Item(22, 2).Value = 7005
$ws.Cells.Item(22, 3).Value = 30731
$ws.Cells.Item(23, 1).Value = 275852
$ws.Cells.Item(23, 2).Value = 11888
$ws.Cells.Item(23, 3).Value = 41129
$ws.Cells.Item(24, 1).Value = 340534
$ws.Cells.Item(24, 2).Value = 10792
$ws.Cells.Item(24, 3).Value = 36047
$ws.Cells.Item(25, 1).Value = 290893
$ws.Cells.Item(25, 2).Value = 10248
$ws.Cells.Item(25, 3).Value = 38098
$ws.Cells.Item(26, 1).Value = 2558656
$ws.Cells.Item(26, 2).Value = 10680
$ws.Cells.Item(26, 3).Value = 36637
$ws.Cells.Item(27, 1).Value = 25327
$ws.Cells.Item(27, 2).Value = 7579
$ws.Cells.Item(27, 3).Value = 44740
$ws.Cells.Item(28, 1).Value = 221751
$ws.Cells.Item(28, 2).Value = 7268
$ws.Cells.Item(28, 3).Value = 36734
$ws.Cells.Item(29, 1).Value = 228116
$ws.Cells.Item(29, 2).Value = 5308
$ws.Cells.Item(29, 3).Value = 33062
$ws.Cells.Item(30, 1).Value = 65024
$ws.Cells.Item(30, 2).Value = 7674
$ws.Cells.Item(30, 3).Value = 43213
$ws.Cells.Item(31, 1).Value = 128837
$ws.Cells.Item(31, 2).Value = 32051
$ws.Cells.Item(31, 3).Value = 39192
$ws.Cells.Item(32, 1).Value = 221337
$ws.Cells.Item(32, 2).Value = 7271
$ws.Cells.Item(32, 3).Value = 42157
$ws.Cells.Item(33, 1).Value = 181460
$ws.Cells.Item(33, 2).Value = 5193
$ws.Cells.Item(33, 3).Value = 22719
$ws.Cells.Item(34, 1).Value = 223208
$ws.Cells.Item(34, 2).Value = 2713
$ws.Cells.Item(34, 3).Value = 24659
$ws.Cells.Item(35, 1).Value = 217685
$ws.Cells.Item(35, 2).Value = 3330
$ws.Cells.Item(35, 3).Value = 23189
$ws.Cells.Item(36, 1).Value = 118492
$ws.Cells.Item(36, 2).Value = 5776
$ws.Cells.Item(36, 3).Value = 26158
$ws.Cells.Item(37, 1).Value = 322886
$ws.Cells.Item(37, 2).Value = 2667
$ws.Cells.Item(37, 3).Value = 16213
$ws.Cells.Item(38, 1).Value = 51454
$ws.Cells.Item(38, 2).Value = 5667
$ws.Cells.Item(38, 3).Value = 4542177
$ws.Cells.Item(39, 1).Value = 17021
$ws.Cells.Item(39, 2).Value = 5800
$ws.Cells.Item(39, 3).Value = 53632
$ws.Cells.Item(40, 1).Value = 190833
$ws.Cells.Item(40, 2).Value = 2319
$ws.Cells.Item(40, 3).Value = 37048
$ws.Cells.Item(41, 1).Value = 197787
$ws.Cells.Item(41, 2).Value = 2395
$ws.Cells.Item(41, 3).Value = 10931
$ws.Cells.Item(42, 1).Value = 172445
$ws.Cells.Item(42, 2).Value = 4036
$ws.Cells.Item(42, 3).Value = 12144
$ws.Cells.Item(43, 1).Value = 50222
$ws.Cells.Item(43, 2).Value = 5036
$ws.Cells.Item(43, 3).Value = 12392
$ws.Cells.Item(44, 1).Value = 111734
$ws.Cells.Item(44, 2).Value = 5752
$ws.Cells.Item(44, 3).Value = 14556
$ws.Cells.Item(45, 1).Value = 32363
$ws.Cells.Item(45, 2).Value = 5456
$ws.Cells.Item(45, 3).Value = 11383
$ws.Cells.Item(46, 1).Value = 200967
$ws.Cells.Item(46, 2).Value = 1351
$ws.Cells.Item(46, 3).Value = 11089
$ws.Cells.Item(47, 1).Value = 161393
$ws.Cells.Item(47, 2).Value = 4468
$ws.Cells.Item(47, 3).Value = 10652
$ws.Cells.Item(48, 1).Value = 264034
$ws.Cells.Item(48, 2).Value = 3332
$ws.Cells.Item(48, 3).Value = 19177
$ws.Cells.Item(49, 1).Value = 155328
$ws.Cells.Item(49, 2).Value = 4711
$ws.Cells.Item(49, 3).Value = 13599
$ws.Cells.Item(50, 1).Value = 166897
$ws.Cells.Item(50, 2).Value = 4363
$ws.Cells.Item(50, 3).Value = 11643
$ws.Cells.Item(51, 1).Value = 17980
$ws.Cells.Item(51, 2).Value = 5515
$ws.Cells.Item(51, 3).Value = 11277
$ws.Cells.Item(52, 1).Value = 172040
$ws.Cells.Item(52, 2).Value = 1240
$ws.Cells.Item(52, 3).Value = 12485
$ws.Cells.Item(53, 1).Value = 121472
$ws.Cells.Item(53, 2).Value = 4646
$ws.Cells.Item(53, 3).Value = 8099
$ws.Cells.Item(54, 1).Value = 85330
$ws.Cells.Item(54, 2).Value = 6612
$ws.Cells.Item(54, 3).Value = 8839
$ws.Cells.Item(55, 1).Value = 64292
$ws.Cells.Item(55, 2).Value = 6379
$ws.Cells.Item(55, 3).Value = 9531
$ws.Cells.Item(56, 1).Value = 142770
$ws.Cells.Item(56, 2).Value = 5141
$ws.Cells.Item(56, 3).Value = 10108
$ws.Cells.Item(57, 1).Value = 43909
$ws.Cells.Item(57, 2).Value = 5917
$ws.Cells.Item(57, 3).Value = 7474
$ws.Cells.Item(58, 1).Value = 21183
$ws.Cells.Item(58, 2).Value = 5804
$ws.Cells.Item(58, 3).Value = 6812
$ws.Cells.Item(59, 1).Value = 24135
$ws.Cells.Item(59, 2).Value = 4480
$ws.Cells.Item(59, 3).Value = 7909
$ws.Cells.Item(60, 1).Value = 98724
$ws.Cells.Item(60, 2).Value = 4315
$ws.Cells.Item(60, 3).Value = 8132
$ws.Cells.Item(61, 1).Value = 187953
$ws.Cells.Item(61, 2).Value = 2316
$ws.Cells.Item(61, 3).Value = 7182
$ws.Cells.Item(62, 1).Value = 51959
$ws.Cells.Item(62, 2).Value = 4870
$ws.Cells.Item(62, 3).Value = 7641
$ws.Cells.Item(63, 1).Value = 129239
$ws.Cells.Item(63, 2).Value = 5015
$ws.Cells.Item(63, 3).Value = 9401
$ws.Cells.Item(64, 1).Value = 94874
$ws.Cells.Item(64, 2).Value = 5606
$ws.Cells.Item(64, 3).Value = 10677
$ws.Cells.Item(65, 1).Value = 162275
$ws.Cells.Item(65, 2).Value = 3486
$ws.Cells.Item(65, 3).Value = 7861
$ws.Cells.Item(66, 1).Value = 28505
$ws.Cells.Item(66, 2).Value = 5157
$ws.Cells.Item(66, 3).Value = 7892
$ws.Cells.Item(67, 1).Value = 149651
$ws.Cells.Item(67, 2).Value = 4246
$ws.Cells.Item(67, 3).Value = 7701
$ws.Cells.Item(68, 1).Value = 13430
$ws.Cells.Item(68, 2).Value = 5706
$ws.Cells.Item(68, 3).Value = 7479
$ws.Cells.Item(69, 1).Value = 7280
$ws.Cells.Item(69, 2).Value = 4981
$ws.Cells.Item(69, 3).Value = 7117
$ws.Cells.Item(70, 1).Value = 180055
$ws.Cells.Item(70, 2).Value = 979
$ws.Cells.Item(70, 3).Value = 7740
$ws.Cells.Item(71, 1).Value = 188406
$ws.Cells.Item(71, 2).Value = 2231
$ws.Cells.Item(71, 3).Value = 10415
$ws.Cells.Item(72, 1).Value = 186513
$ws.Cells.Item(72, 2).Value = 3542
$ws.Cells.Item(72, 3).Value = 8027
$ws.Cells.Item(73, 1).Value = 29356
$ws.Cells.Item(73, 2).Value = 5127
$ws.Cells.Item(73, 3).Value = 10498
$ws.Cells.Item(74, 1).Value = 80032
$ws.Cells.Item(74, 2).Value = 5436
$ws.Cells.Item(74, 3).Value = 7722
$ws.Cells.Item(75, 1).Value = 101055
$ws.Cells.Item(75, 2).Value = 5670
$ws.Cells.Item(75, 3).Value = 8486
$ws.Cells.Item(76, 1).Value = 187518
$ws.Cells.Item(76, 2).Value = 663
$ws.Cells.Item(76, 3).Value = 6135
$ws.Cells.Item(77, 1).Value = 7502
$ws.Cells.Item(77, 2).Value = 4387
$ws.Cells.Item(77, 3).Value = 7304
$ws.Cells.Item(78, 1).Value = 171799
$ws.Cells.Item(78, 2).Value = 5199
$ws.Cells.Item(78, 3).Value = 7621
$ws.Cells.Item(79, 1).Value = 184277
$ws.Cells.Item(79, 2).Value = 3203
$ws.Cells.Item(79, 3).Value = 8744
$ws.Cells.Item(80, 1).Value = 7135
$ws.Cells.Item(80, 2).Value = 5636
$ws.Cells.Item(80, 3).Value = 6804
$ws.Cells.Item(81, 1).Value = 57426
$ws.Cells.Item(81, 2).Value = 5128
$ws.Cells.Item(81, 3).Value = 9005
$ws.Cells.Item(82, 1).Value = 176634
$ws.Cells.Item(82, 2).Value = 1885
$ws.Cells.Item(82, 3).Value = 7320
$ws.Cells.Item(83, 1).Value = 19681
$ws.Cells.Item(83, 2).Value = 5660
$ws.Cells.Item(83, 3).Value = 6827
$ws.Cells.Item(84, 1).Value = 7616
$ws.Cells.Item(84, 2).Value = 5931
$ws.Cells.Item(84, 3).Value = 6761
$ws.Cells.Item(85, 1).Value = 192737
$ws.Cells.Item(85, 2).Value = 1401
$ws.Cells.Item(85, 3).Value = 8362
$ws.Cells.Item(86, 1).Value = 166321
$ws.Cells.Item(86, 2).Value = 3094
$ws.Cells.Item(86, 3).Value = 7700
$ws.Cells.Item(87, 1).Value = 108373
$ws.Cells.Item(87, 2).Value = 4655
$ws.Cells.Item(87, 3).Value = 7828
$ws.Cells.Item(88, 1).Value = 160086
$ws.Cells.Item(88, 2).Value = 4327
$ws.Cells.Item(88, 3).Value = 8221
$ws.Cells.Item(89, 1).Value = 34978
$ws.Cells.Item(89, 2).Value = 5086
$ws.Cells.Item(89, 3).Value = 7461
$ws.Cells.Item(90, 1).Value = 121946
$ws.Cells.Item(90, 2).Value = 3681
$ws.Cells.Item(90, 3).Value = 8827
$ws.Cells.Item(91, 1).Value = 155497
$ws.Cells.Item(91, 2).Value = 4297
$ws.Cells.Item(91, 3).Value = 8173
$ws.Cells.Item(92, 1).Value = 44685
$ws.Cells.Item(92, 2).Value = 5847
$ws.Cells.Item(92, 3).Value = 7986
$ws.Cells.Item(93, 1).Value = 97633
$ws.Cells.Item(93, 2).Value = 5379
$ws.Cells.Item(93, 3).Value = 8259
$ws.Cells.Item(94, 1).Value = 42531
$ws.Cells.Item(94, 2).Value = 6178
$ws.Cells.Item(94, 3).Value = 8549
$ws.Cells.Item(95, 1).Value = 10655
$ws.Cells.Item(95, 2).Value = 6201
$ws.Cells.Item(95, 3).Value = 7398
$ws.Cells.Item(96, 1).Value = 170209
$ws.Cells.Item(96, 2).Value = 3768
$ws.Cells.Item(96, 3).Value = 8829
$ws.Cells.Item(97, 1).Value = 182474
$ws.Cells.Item(97, 2).Value = 3018
$ws.Cells.Item(97, 3).Value = 7962
$ws.Cells.Item(98, 1).Value = 72093
$ws.Cells.Item(98, 2).Value = 4556
$ws.Cells.Item(98, 3).Value = 15780
$ws.Cells.Item(99, 1).Value = 62959
$ws.Cells.Item(99, 2).Value = 4754
$ws.Cells.Item(99, 3).Value = 11393
$ws.Cells.Item(100, 1).Value = 6892
$ws.Cells.Item(100, 2).Value = 4736
$ws.Cells.Item(100, 3).Value = 7053
$ws.Cells.Item(101, 1).Value = 34654
$ws.Cells.Item(101, 2).Value = 5734
$ws.Cells.Item(101, 3).Value = 7513
$ws.Cells.Item(102, 1).Value = 178638
$ws.Cells.Item(102, 2).Value = 1217
$ws.Cells.Item(102, 3).Value = 5869
